# Auto-generated edit script: updates Price (D) and Volume(1h) (E) columns
# to reflect the latest cryptos snapshot, per commit diff.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '66.547.10'
$ws.Range("E2").Value = '  +0.69%  '
$ws.Range("D3").Value = '3.230.57'
$ws.Range("E3").Value = '  +1.59%  '
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '0.999'
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = '  -0.13%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '605.33'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +1.81%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '157.83'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  +3.02%  '
$ws.Range("E7").Value = '  -0.05%  '
$ws.Range("D8").Value = '3.229.34'
$ws.Range("E8").Value = '  +1.62%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.547'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  +0.92%  '
$ws.Range("E10").Value = '  +1.73%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '5.68'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  -5.44%  '
$ws.Range("E12").Value = '  -1.51%  '
$ws.Range("E13").Value = '  +3.08%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '38.92'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  +0.19%  '
$ws.Range("D15").Value = '3.760.79'
$ws.Range("E15").Value = '  +1.58%  '
$ws.Range("D16").Value = '66.625.46'
$ws.Range("E16").Value = '  +0.81%  '
$ws.Range("E17").Value = '  +0.07%  '
$ws.Range("D18").Value = '3.228.78'
$ws.Range("E18").Value = '  +1.36%  '
$ws.Range("E19").Value = '  +1.12%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '510.34'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  +0.65%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '15.23'
$ws.Range("D21").Style = "Normal"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '0.734'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  +0.10%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '8.05'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  +0.79%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '14.69'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  -2.26%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '84.89'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  +0.31%  '
$ws.Range("E26").Value = '  +0.27%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '2.99'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  +0.13%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '9.17'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  +0.19%  '
$ws.Range("E29").Value = '  +4.65%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '2.96'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  +3.29%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '7.02'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  +0.30%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '28.20'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  +0.61%  '
$ws.Range("E33").Value = '  +0.01%  '
$ws.Range("E34").Value = '  -2.79%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '0.105'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  +17.65%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '6.50'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  +0.48%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '508.15'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  +5.00%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '55.64'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  +1.79%  '
$ws.Range("E39").Value = '  +17.86%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.0423'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  +0.82%  '
$ws.Range("E41").Value = '  +7.80%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.129'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  +6.59%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '8.74'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  -1.29%  '
$ws.Range("E44").Value = '  +0.20%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '2.46'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  +2.56%  '
$ws.Range("D46").Value = '2.872.23'
$ws.Range("E46").Value = '  -0.85%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '28.53'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  +0.89%  '
$ws.Range("E48").Value = '  +4.68%  '
$ws.Range("E50").Value = '  -0.14%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '122.26'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  +0.51%  '
